# Add rain-barrel storage parameters and fix simulate_performances() /
# control.curve row on the "rain_barrel" sheet, then leave that sheet as
# the active / selected tab (mirrors the author's manual edit session).

$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("rain_barrel")

# ---------------------------------------------------------------------
# 1) Three new "storage" parameter rows (id 2,3,4) pushed in above the
#    existing "drain" rows. New shared strings are written in this exact
#    order so they land at uniqueCount indices 108-111, matching the
#    target workbook (unknown1, unknown2, comment, unknown3).
# ---------------------------------------------------------------------
$comment = "add parameter generated automatically in GUI, otherwise run fails"

# row 4 : storage / id 2 / unknown1
$ws4.Cells.Item(4, 1).Value2 = "storage"
$ws4.Cells.Item(4, 2).Value2 = 2
$ws4.Cells.Item(4, 3).Value2 = "unknown1"
$ws4.Cells.Item(4, 4).Value2 = 0.75
$ws4.Cells.Item(4, 5).Value2 = 0.75
$ws4.Cells.Item(4, 6).Value2 = 0.75
$ws4.Cells.Item(4, 7).Value2 = 0.75
$ws4.Cells.Item(4, 8).Value2 = $comment

# row 5 : storage / id 3 / unknown2
$ws4.Cells.Item(5, 1).Value2 = "storage"
$ws4.Cells.Item(5, 2).Value2 = 3
$ws4.Cells.Item(5, 3).Value2 = "unknown2"
$ws4.Cells.Item(5, 4).Value2 = 0.5
$ws4.Cells.Item(5, 5).Value2 = 0.5
$ws4.Cells.Item(5, 6).Value2 = 0.5
$ws4.Cells.Item(5, 7).Value2 = 0.5
$ws4.Cells.Item(5, 8).Value2 = $comment

# row 6 : storage / id 4 / unknown3
$ws4.Cells.Item(6, 1).Value2 = "storage"
$ws4.Cells.Item(6, 2).Value2 = 4
$ws4.Cells.Item(6, 3).Value2 = "unknown3"
$ws4.Cells.Item(6, 4).Value2 = 0
$ws4.Cells.Item(6, 5).Value2 = 0
$ws4.Cells.Item(6, 6).Value2 = 0
$ws4.Cells.Item(6, 7).Value2 = 0
$ws4.Cells.Item(6, 8).Value2 = $comment

# ---------------------------------------------------------------------
# 2) The pre-existing "drain" parameter rows, shifted down by 3 rows
#    (old row 4 -> new row 7, ... old row 9 -> new row 12). Values are
#    unchanged except the last ("control.curve") row, which previously
#    had only a blank placeholder cell and now is fully populated
#    (fix for simulate_performances()).
# ---------------------------------------------------------------------

# row 7 : drain / id 1 / flow.coefficient  (unchanged content)
$ws4.Cells.Item(7, 1).Value2 = "drain"
$ws4.Cells.Item(7, 2).Value2 = 1
$ws4.Cells.Item(7, 3).Value2 = "flow.coefficient"
$ws4.Cells.Item(7, 4).Value2 = 4
$ws4.Cells.Item(7, 5).Value2 = 4
$ws4.Cells.Item(7, 6).Value2 = 4
$ws4.Cells.Item(7, 7).Value2 = 4
$ws4.Cells.Item(7, 8).Value2 = 4

# row 8 : drain / id 2 / flow.exponent  (unchanged content)
$ws4.Cells.Item(8, 1).Value2 = "drain"
$ws4.Cells.Item(8, 2).Value2 = 2
$ws4.Cells.Item(8, 3).Value2 = "flow.exponent"
$ws4.Cells.Item(8, 4).Value2 = 0.5
$ws4.Cells.Item(8, 5).Value2 = 0.5
$ws4.Cells.Item(8, 6).Value2 = 0.5
$ws4.Cells.Item(8, 7).Value2 = 0.5
$ws4.Cells.Item(8, 8).Value2 = 0.5

# row 9 : drain / id 3 / offset_mm  (unchanged content)
$ws4.Cells.Item(9, 1).Value2 = "drain"
$ws4.Cells.Item(9, 2).Value2 = 3
$ws4.Cells.Item(9, 3).Value2 = "offset_mm"
$ws4.Cells.Item(9, 4).Value2 = 0
$ws4.Cells.Item(9, 5).Value2 = 200
$ws4.Cells.Item(9, 6).Value2 = 400
$ws4.Cells.Item(9, 7).Value2 = 600

# row 10 : drain / id 4 / level.open_mm  (unchanged content)
$ws4.Cells.Item(10, 1).Value2 = "drain"
$ws4.Cells.Item(10, 2).Value2 = 4
$ws4.Cells.Item(10, 3).Value2 = "level.open_mm"
$ws4.Cells.Item(10, 4).Value2 = 0.1
$ws4.Cells.Item(10, 5).Value2 = 0.1
$ws4.Cells.Item(10, 6).Value2 = 0.1
$ws4.Cells.Item(10, 7).Value2 = 0.1
$ws4.Cells.Item(10, 8).Value2 = 0

# row 11 : drain / id 5 / level.closed_mm  (unchanged content)
$ws4.Cells.Item(11, 1).Value2 = "drain"
$ws4.Cells.Item(11, 2).Value2 = 5
$ws4.Cells.Item(11, 3).Value2 = "level.closed_mm"
$ws4.Cells.Item(11, 4).Value2 = 0
$ws4.Cells.Item(11, 5).Value2 = 0
$ws4.Cells.Item(11, 6).Value2 = 0
$ws4.Cells.Item(11, 7).Value2 = 0
$ws4.Cells.Item(11, 8).Value2 = 0

# row 12 : drain / id 6 / control.curve  (was an empty placeholder row;
# now filled in so simulate_performances() no longer fails)
$ws4.Cells.Item(12, 1).Value2 = "drain"
$ws4.Cells.Item(12, 2).Value2 = 6
$ws4.Cells.Item(12, 3).Value2 = "control.curve"
$ws4.Cells.Item(12, 4).Value2 = 0
$ws4.Cells.Item(12, 5).Value2 = 0
$ws4.Cells.Item(12, 6).Value2 = 0
$ws4.Cells.Item(12, 7).Value2 = 0
$ws4.Cells.Item(12, 8).Value2 = $comment

# ---------------------------------------------------------------------
# 3) Give the "storage"/"drain" columns A-C a thin box border on every
#    cell in the block (matches the look of the rest of the table) -
#    applied per-cell so each one gets its own full border, not just the
#    outer edge of the range.
# ---------------------------------------------------------------------
foreach ($r in 4..12) {
    foreach ($c in 1..3) {
        $cell = $ws4.Cells.Item($r, $c)
        $cell.Borders.Item(7).LineStyle  = 1
        $cell.Borders.Item(7).Weight     = 2
        $cell.Borders.Item(8).LineStyle  = 1
        $cell.Borders.Item(8).Weight     = 2
        $cell.Borders.Item(9).LineStyle  = 1
        $cell.Borders.Item(9).Weight     = 2
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight    = 2
    }
    $ws4.Cells.Item($r, 1).Font.Bold = $true
    $ws4.Cells.Item($r, 2).Font.Bold = $true
}

# ---------------------------------------------------------------------
# 4) View state: rain_barrel becomes the active / selected sheet
#    (green_roof loses tabSelected automatically), scrolled/selected on
#    its last edited cell.
# ---------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("H12").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
